$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, copying the formatting of the existing header cells
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
